$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 86419294
$ws.Range("B4").Value = 77506
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("Q4").Value = 440686.1394479795
$ws.Range("R4").Value = 6707147.171128325
$ws.Range("AC4").Value = "På tall"

# Row 5
$ws.Range("A5").Value = 86419305
$ws.Range("B5").Value = 77506
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").Value = "(Ach.) Ach."
$ws.Range("Q5").Value = 440606.8734944779
$ws.Range("R5").Value = 6707280.052989913
$ws.Range("AC5").Value = "På flera tallar"

# Row 6
$ws.Range("A6").Value = 86419313
$ws.Range("B6").Value = 77506
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("Q6").Value = 440607.1726549119
$ws.Range("R6").Value = 6707238.159541409
$ws.Range("AC6").Value = "rikligt på flera tallar"

# Row 7
$ws.Range("A7").Value = 86419293
$ws.Range("B7").Value = 77506
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("Q7").Value = 440687.1425972193
$ws.Range("R7").Value = 6707148.140317255
$ws.Range("AC7").ClearContents()

# Row 8
$ws.Range("A8").Value = 86419304
$ws.Range("B8").Value = 90653
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 4364
$ws.Range("F8").Value = "Dropptaggsvamp"
$ws.Range("G8").Value = "Hydnellum ferrugineum"
$ws.Range("H8").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q8").Value = 440783.8270494898
$ws.Range("R8").Value = 6707144.091754919
$ws.Range("AC8").ClearContents()

# Row 9
$ws.Range("A9").Value = 86419296
$ws.Range("B9").Value = 77506
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 6425
$ws.Range("F9").Value = "Garnlav"
$ws.Range("G9").Value = "Alectoria sarmentosa"
$ws.Range("H9").Value = "(Ach.) Ach."
$ws.Range("Q9").Value = 440814.1817916233
$ws.Range("R9").Value = 6707128.810180089
$ws.Range("AC9").Value = "På tall"

# Row 10
$ws.Range("A10").Value = 86419290
$ws.Range("B10").Value = 8377
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 106545
$ws.Range("F10").Value = "Mindre märgborre"
$ws.Range("G10").Value = "Tomicus minor"
$ws.Range("H10").Value = "(Hartig, 1834)"
$ws.Range("Q10").Value = 440814.1656648018
$ws.Range("R10").Value = 6707127.824798071
$ws.Range("AC10").ClearContents()
